# feat: fix bad formatting
# Reorders the category rows 27-42 on the "Categories" sheet back into
# their correct (intended) sequence. Only columns B:G change; columns
# A (data-object) and H (parent.typeId) stay "category" throughout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for columns B (key), C (description), D (externalId),
# E (name), F (slug), G (parent.key) for rows 27 through 42, in the
# corrected order.
$rows = @(
    @("child-anyKey",      "child-anyDescription",      "child-anyId",      "Child-any",      "child-anySlug",      "ageGroupKey"),
    @("JonesKey",          "JonesDescription",          "JonesId",          "Jones",          "JonesSlug",          "brandKey"),
    @("BataleonKey",       "BataleonDescription",       "BataleonId",       "Bataleon",       "BataleonSlug",       "brandKey"),
    @("child-boyKey",      "child-boyDescription",      "child-boyId",      "Child-boy",      "child-boySlug",      "ageGroupKey"),
    @("SalomonKey",        "SalomonDescription",        "SalomonId",        "Salomon",        "SalomonSlug",        "brandKey"),
    @("RomeKey",           "RomeDescription",           "RomeId",           "Rome",           "RomeSlug",           "brandKey"),
    @("SimsKey",           "SimsDescription",           "SimsId",           "Sims",           "SimsSlug",           "brandKey"),
    @("SplitboardingKey",  "SplitboardingDescription",  "SplitboardingId",  "Splitboarding",  "SplitboardingSlug",  "terrainKey"),
    @("NitroKey",          "NitroDescription",          "NitroId",          "Nitro",          "NitroSlug",          "brandKey"),
    @("NeverSummerKey",    "NeverSummerDescription",    "NeverSummerId",    "NeverSummer",    "NeverSummerSlug",    "brandKey"),
    @("child-girlKey",     "child-girlDescription",     "child-girlId",     "Child-girl",     "child-girlSlug",     "ageGroupKey"),
    @("UnitedShapesKey",   "UnitedShapesDescription",   "UnitedShapesId",   "UnitedShapes",   "UnitedShapesSlug",   "brandKey"),
    @("3YearsKey",         "3YearsDescription",         "3YearsId",         "3Years",         "3YearsSlug",         "terrainKey"),
    @("YesKey",            "YesDescription",            "YesId",            "Yes",            "YesSlug",            "brandKey"),
    @("WestonKey",         "WestonDescription",         "WestonId",         "Weston",         "WestonSlug",         "brandKey"),
    @("CardiffKey",        "CardiffDescription",        "CardiffId",        "Cardiff",        "CardiffSlug",        "brandKey")
)

$startRow = 27
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 2).Value = $data[0]
    $ws.Cells.Item($r, 3).Value = $data[1]
    $ws.Cells.Item($r, 4).Value = $data[2]
    $ws.Cells.Item($r, 5).Value = $data[3]
    $ws.Cells.Item($r, 6).Value = $data[4]
    $ws.Cells.Item($r, 7).Value = $data[5]
}
